# Applies the "converted calibration pts readings to pd.DF" edit:
#   - On the "instructions" sheet, draw a bordered "box" around the
#     existing block of help text (rows 11-26, ending on the "Example"
#     label which is moved up from B28 to B26), and around a second
#     block (rows 29-32) that gets a brand new closing note appended
#     in B32.
#   - Update the active-cell selection to B22.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("instructions")

# --- Box 1: B11 (bold header) .. B26 (former B28, now the closing line) ---

# Header row B11 is already bold; add the top/left/right edge of the box.
$ws.Range("B11").Borders.Item(7).LineStyle = 1
$ws.Range("B11").Borders.Item(10).LineStyle = 1
$ws.Range("B11").Borders.Item(8).LineStyle = 1

# Middle rows B12:B25 just get the left/right edges of the box.
$midBox1 = $ws.Range("B12:B25")
$midBox1.Borders.Item(7).LineStyle = 1
$midBox1.Borders.Item(10).LineStyle = 1

# Move the "Example" label from B28 up to B26 and make it the (non-bold)
# closing line of the box (left/right/bottom edges).
$exampleValue = $ws.Range("B28").Value2
$ws.Range("B28").Clear()
$ws.Range("B26").Value2 = $exampleValue
$ws.Range("B26").Font.Bold = $false
$ws.Range("B26").Borders.Item(7).LineStyle = 1
$ws.Range("B26").Borders.Item(10).LineStyle = 1
$ws.Range("B26").Borders.Item(9).LineStyle = 1

# --- Box 2: B29 (bold header) .. B32 (new closing note) ---

$ws.Range("B29").Font.Bold = $true
$ws.Range("B29").Borders.Item(7).LineStyle = 1
$ws.Range("B29").Borders.Item(10).LineStyle = 1
$ws.Range("B29").Borders.Item(8).LineStyle = 1

$midBox2 = $ws.Range("B30:B31")
$midBox2.Borders.Item(7).LineStyle = 1
$midBox2.Borders.Item(10).LineStyle = 1

$ws.Range("B32").Value2 = "Important: the order of quantities must correspond to the same order defined for ‘Calibration quantity i' in the user input."
$ws.Range("B32").Font.Bold = $false
$ws.Range("B32").Borders.Item(7).LineStyle = 1
$ws.Range("B32").Borders.Item(10).LineStyle = 1
$ws.Range("B32").Borders.Item(9).LineStyle = 1

# Match the author's final cursor position.
$ws.Range("B22").Select()
